# =============================================================================
# feat: add 2022-Q1 data
#
# The workbook's old last sheet "总计" (the running roll-up) is repurposed as
# the new "2022-Q1" detail sheet (it keeps its original sheetPr/styles), and a
# fresh copy of it becomes the new "总计" roll-up appended after it — this
# mirrors the sheetId/rId renumbering in the target diff (2022-Q1 <- old
# sheetId 6, 总计 <- new sheetId 7) and keeps every sheet's boilerplate
# (outline/page-setup props) intact.
# =============================================================================

$wb = $excel.ActiveWorkbook

$q1Sheet = $wb.Worksheets.Item("总计")
$q1Sheet.Name = "2022-Q1"

$q1Sheet.Copy($null, $q1Sheet)
$totalSheet = $wb.Worksheets.Item("2022-Q1 (2)")
$totalSheet.Name = "总计"

# ---------------------------------------------------------------------------
# 1. "2022-Q1" fund-holding detail sheet.
# ---------------------------------------------------------------------------

# Extend the existing bold/centered-header style (column A index + header
# row) from the old 4-column/6-row layout to the new 8-column/19-row one,
# reusing the workbook's existing style index instead of minting new ones.
$q1Sheet.Range("A6").Copy()
$q1Sheet.Range("A7:A19").PasteSpecial(-4122)
$q1Sheet.Range("D1").Copy()
$q1Sheet.Range("E1:H1").PasteSpecial(-4122)

$q1Sheet.Range("B1").Value = "基金代码"
$q1Sheet.Range("C1").Value = "基金名称"
$q1Sheet.Range("D1").Value = "基金规模"
$q1Sheet.Range("E1").Value = "股票总仓位"
$q1Sheet.Range("F1").Value = "仓位占比"
$q1Sheet.Range("G1").Value = "持有市值(亿元)"
$q1Sheet.Range("H1").Value = "仓位排名"

# Columns B-G are text-typed (fund code/name/numbers stored as text, matching
# the source data); A is the row index, H is the numeric position rank.
$fundRows = @(
    @(0, "515210", "国泰中证钢铁ETF", "16.24", "99.25", "5.32", "0.8640", 5),
    @(1, "502023", "鹏华国证钢铁行业指数（LOF）", "15.55", "94.76", "5.32", "0.8273", 4),
    @(2, "005273", "华商可转债债券A", "10.06", "39.60", "2.79", "0.2807", 3),
    @(3, "168203", "中融国证钢铁行业指数", "4.30", "92.58", "5.17", "0.2223", 4),
    @(4, "004495", "博时量化平衡混合", "10.49", "38.32", "1.17", "0.1227", 5),
    @(5, "005284", "华商可转债债券C", "3.80", "39.60", "2.79", "0.1060", 3),
    @(6, "002317", "招商睿逸稳健配置混合", "6.08", "48.74", "1.34", "0.0815", 10),
    @(7, "002311", "创金合信中证500指数增强A", "5.72", "92.95", "1.10", "0.0629", 5),
    @(8, "501059", "西部利得中证国有企业红利指数增强（LOF）A", "2.22", "92.96", "2.81", "0.0624", 9),
    @(9, "012963", "招商稳健平衡混合A", "1.33", "64.08", "2.22", "0.0295", 9),
    @(10, "002316", "创金合信中证500指数增强C", "2.60", "92.95", "1.10", "0.0286", 5),
    @(11, "009439", "西部利得中证国有企业红利指数增强（LOF）C", "0.98", "92.96", "2.81", "0.0275", 9),
    @(12, "012964", "招商稳健平衡混合C", "0.59", "64.08", "2.22", "0.0131", 9),
    @(13, "260117", "景顺长城支柱产业混合", "0.24", "76.78", "3.04", "0.0073", 7),
    @(14, "013802", "财通资管中证钢铁指数A", "0.11", "90.83", "5.13", "0.0056", 3),
    @(15, "008838", "德邦量化对冲策略灵活配置混合A", "0.20", "38.81", "1.19", "0.0024", 8),
    @(16, "013803", "财通资管中证钢铁指数C", "0.02", "90.83", "5.13", "0.0010", 3),
    @(17, "008839", "德邦量化对冲策略灵活配置混合C", "0.02", "38.81", "1.19", "0.0002", 8)
)

$r = 2
foreach ($row in $fundRows) {
    $q1Sheet.Cells.Item($r, 1).Value = $row[0]

    for ($c = 2; $c -le 7; $c++) {
        $cell = $q1Sheet.Cells.Item($r, $c)
        $cell.NumberFormat = "@"
        $cell.Value = $row[$c - 1]
        $cell.Style = "Normal"
    }

    $q1Sheet.Cells.Item($r, 8).Value = $row[7]
    $r = $r + 1
}

# ---------------------------------------------------------------------------
# 2. "总计" roll-up sheet: prepend the 2022-Q1 summary row and renumber the
#    leading index column (A) so it stays sequential.
# ---------------------------------------------------------------------------
$totalSheet.Range("A6").Copy()
$totalSheet.Range("A7").PasteSpecial(-4122)

$summaryRows = @(
    @(0, "2022-Q1", 18, 2.74),
    @(1, "2021-Q4", 19, 2.69),
    @(2, "2021-Q3", 29, 7.14),
    @(3, "2021-Q2", 22, 4.17),
    @(4, "2021-Q1", 37, 4.15),
    @(5, "2020-Q4", 16, 0.53)
)

$r = 2
foreach ($row in $summaryRows) {
    $totalSheet.Cells.Item($r, 1).Value = $row[0]
    $totalSheet.Cells.Item($r, 2).Value = $row[1]
    $totalSheet.Cells.Item($r, 3).Value = $row[2]
    $totalSheet.Cells.Item($r, 4).Value = $row[3]
    $r = $r + 1
}
